$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-20 Saturday", 2) | Out-Null
$d.Content.Find.Execute("41×91=3731", $true, $false, $false, $false, $false, $true, 1, $false, "42×60=2520", 2) | Out-Null
$d.Content.Find.Execute("58×50=2900", $true, $false, $false, $false, $false, $true, 1, $false, "89×82=7298", 2) | Out-Null
$d.Content.Find.Execute("28×95=2660", $true, $false, $false, $false, $false, $true, 1, $false, "93×78=7254", 2) | Out-Null
$d.Content.Find.Execute("12×12=144", $true, $false, $false, $false, $false, $true, 1, $false, "19×89=1691", 2) | Out-Null
$d.Content.Find.Execute("33×21=693", $true, $false, $false, $false, $false, $true, 1, $false, "48×20=960", 2) | Out-Null
$d.Content.Find.Execute("31×14=434", $true, $false, $false, $false, $false, $true, 1, $false, "35×86=3010", 2) | Out-Null
$d.Content.Find.Execute("59×69=4071", $true, $false, $false, $false, $false, $true, 1, $false, "68×18=1224", 2) | Out-Null
$d.Content.Find.Execute("41×38=1558", $true, $false, $false, $false, $false, $true, 1, $false, "49×65=3185", 2) | Out-Null
$d.Content.Find.Execute("74×53=3922", $true, $false, $false, $false, $false, $true, 1, $false, "88×63=5544", 2) | Out-Null
$d.Content.Find.Execute("33×59=1947", $true, $false, $false, $false, $false, $true, 1, $false, "97×65=6305", 2) | Out-Null
$d.Content.Find.Execute("28×75=2100", $true, $false, $false, $false, $false, $true, 1, $false, "33×30=990", 2) | Out-Null
$d.Content.Find.Execute("43×71=3053", $true, $false, $false, $false, $false, $true, 1, $false, "98×81=7938", 2) | Out-Null
$d.Content.Find.Execute("90×30=2700", $true, $false, $false, $false, $false, $true, 1, $false, "11×97=1067", 2) | Out-Null
$d.Content.Find.Execute("53×98=5194", $true, $false, $false, $false, $false, $true, 1, $false, "98×90=8820", 2) | Out-Null
$d.Content.Find.Execute("11×35=385", $true, $false, $false, $false, $false, $true, 1, $false, "92×94=8648", 2) | Out-Null
$d.Content.Find.Execute("23×22=506", $true, $false, $false, $false, $false, $true, 1, $false, "28×49=1372", 2) | Out-Null
$d.Content.Find.Execute("20×37=740", $true, $false, $false, $false, $false, $true, 1, $false, "54×61=3294", 2) | Out-Null
$d.Content.Find.Execute("82×40=3280", $true, $false, $false, $false, $false, $true, 1, $false, "17×38=646", 2) | Out-Null
$d.Content.Find.Execute("85×84=7140", $true, $false, $false, $false, $false, $true, 1, $false, "75×58=4350", 2) | Out-Null
$d.Content.Find.Execute("25×74=1850", $true, $false, $false, $false, $false, $true, 1, $false, "24×53=1272", 2) | Out-Null
$d.Content.Find.Execute("63×27=1701", $true, $false, $false, $false, $false, $true, 1, $false, "34×77=2618", 2) | Out-Null
$d.Content.Find.Execute("57×87=4959", $true, $false, $false, $false, $false, $true, 1, $false, "24×42=1008", 2) | Out-Null
$d.Content.Find.Execute("87×20=1740", $true, $false, $false, $false, $false, $true, 1, $false, "37×33=1221", 2) | Out-Null
$d.Content.Find.Execute("72×27=1944", $true, $false, $false, $false, $false, $true, 1, $false, "14×75=1050", 2) | Out-Null
$d.Content.Find.Execute("18×56=1008", $true, $false, $false, $false, $false, $true, 1, $false, "50×20=1000", 2) | Out-Null
